$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a few existing timing values (rows 5-7) ---
$ws.Range("B5").Value = 0.0006268024444580078
$ws.Range("B6").Value = 0.0007231235504150391
$ws.Range("B7").Value = 0.0008211135864257812

# --- A8 / A44: tuple-printed lists rewritten as proper JSON-style lists ---
$ws.Range("A8").Value = "[[0, 2], [0, 0], [2, 1], [1, 2], [1, 0], [1, 1], [0, 1]]"
$ws.Range("A44").Value = "[[0, 2], [0, 1], [1, 2], [0, 0], [1, 1], [1, 0], [2, 2]]"

# --- Insert a new row at 68 for "move_fidelity" (everything below shifts down by one) ---
$ws.Rows.Item(68).Insert()

$ws.Range("A68").Value = "move_fidelity"
$ws.Range("B68").Value = 0.9997046993085594

# --- The former row 71 ("total time:") is now row 72; update its value ---
$ws.Range("B72").Value = 0.009095191955566406

Write-Host "done"
